$d = $word.ActiveDocument

# --- (1) & (2) under 5223.370-5 Contract clauses: replace the directly
#     applied tab stop with the "List 2" paragraph style.
$p9 = $d.Paragraphs.Item(9)
$p9.Style = "List 2"

$p10 = $d.Paragraphs.Item(10)
$p10.Style = "List 2"

# --- "(e) Submit requests for waiver ..." paragraph: replace the directly
#     applied bold/black paragraph-mark formatting with the "List 1" style,
#     while preserving the single bold run that sits mid-paragraph.
$p23 = $d.Paragraphs.Item(23)
$boldRunStart = $p23.Range.Start + 78
$boldRun = $d.Range($boldRunStart, $boldRunStart + 1)

$p23.Range.Select()
$word.Selection.ClearFormatting()
$p23.Style = "List 1"
$boldRun.Font.Bold = -1
